$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value  = -0.2769
$ws.Range("B4").Value  = -0.007
$ws.Range("B5").Value  = 0.0346
$ws.Range("B6").Value  = -0.196
$ws.Range("B7").Value  = -0.3249
$ws.Range("B8").Value  = 0.0121
$ws.Range("B9").Value  = -0.0603
$ws.Range("B10").Value = -0.011
$ws.Range("B11").Value = 0.0124
$ws.Range("B12").Value = -0.4636
$ws.Range("B13").Value = -0.0081
$ws.Range("B14").Value = -0.3378
$ws.Range("B15").Value = -0.1039
$ws.Range("B16").Value = -0.0549
$ws.Range("B17").Value = -0.0001
$ws.Range("B18").Value = -0.1122
$ws.Range("B19").Value = 0.0257
$ws.Range("B20").Value = -0.0001
$ws.Range("B21").Value = 0.0131
$ws.Range("B22").Value = -0.0005999999999999999
$ws.Range("B23").Value = -0.0159
$ws.Range("B24").Value = -0.0501
